# Adds the newest "known error" log entry harvested from the local install
# log (keywords column etc.) as the next data row on Sheet1, following the
# same pattern as the existing rows above it.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("A29").Value = "No property enclosed in raw argument brackets: []"
$ws.Range("B29").Value = "10:31:09.010 Error: No property enclosed in raw argument brackets: []"
$ws.Range("C29").Value = "This has no effect on product"
$ws.Range("D29").Value = "Risa"
$ws.Range("E29").Value = "8/25/2014"
$ws.Range("E29").NumberFormat = "d\-mmm\-yy"
$ws.Range("F29").Value = "miniBIP"
$ws.Range("H29").Value = "Linux"
$ws.Range("I29").Value = "Aurora_40_SP_REL"
$ws.Range("J29").Value = "1475_greatest"

$ws.Range("B29:C29").WrapText = $true

$ws.Range("D37").Select() | Out-Null
